$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (preserving things like "1.000" or
# "28.398.13" instead of letting Excel auto-convert them to numbers/dates),
# while leaving the cell's style/format untouched in the saved file.
function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '28.398.13'
$ws.Range("E2").Value = '  -0.24%  '
Set-TextValue "D3" '1.811.28'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  -0.38%  '
Set-TextValue "D5" '313.18'
$ws.Range("E5").Value = '  -0.97%  '
Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  -0.22%  '
Set-TextValue "D7" '0.5167'
$ws.Range("E7").Value = '  -0.38%  '
Set-TextValue "D8" '0.4014'
$ws.Range("E8").Value = '  +3.81%  '
Set-TextValue "D9" '0.07890'
$ws.Range("E9").Value = '  -4.68%  '
Set-TextValue "D10" '1.116'
$ws.Range("E10").Value = '  -0.58%  '
Set-TextValue "D11" '40.94'
$ws.Range("E11").Value = '  -2.17%  '
Set-TextValue "D12" '6.371'
$ws.Range("E12").Value = '  -0.01%  '
Set-TextValue "D13" '1.000'
$ws.Range("E13").Value = '  -0.36%  '
Set-TextValue "D14" '20.45'
$ws.Range("E14").Value = '  -3.22%  '
Set-TextValue "D15" '7.334'
$ws.Range("E15").Value = '  -2.01%  '
Set-TextValue "D16" '1.806.29'
$ws.Range("E16").Value = '  -1.23%  '
Set-TextValue "D17" '92.87'
$ws.Range("E17").Value = '  -1.11%  '
Set-TextValue "D18" '0.00001088'
$ws.Range("E18").Value = '  -2.81%  '
Set-TextValue "D19" '0.06572'
$ws.Range("E19").Value = '  -0.96%  '
Set-TextValue "D20" '1.000'
$ws.Range("E20").Value = '  -0.29%  '
Set-TextValue "D21" '17.35'
$ws.Range("E21").Value = '  -2.66%  '
Set-TextValue "D22" '6.013'
$ws.Range("E22").Value = '  -0.92%  '
Set-TextValue "D23" '28.419.17'
$ws.Range("E23").Value = '  -0.33%  '
Set-TextValue "D24" '11.15'
$ws.Range("E24").Value = '  -3.02%  '
Set-TextValue "D25" '2.229'
$ws.Range("E25").Value = '  -0.65%  '
Set-TextValue "D26" '160.75'
$ws.Range("E26").Value = '  +0.76%  '
Set-TextValue "D27" '20.57'
$ws.Range("E27").Value = '  -2.47%  '
Set-TextValue "D28" '2.023.12'
$ws.Range("E28").Value = '  -0.73%  '
Set-TextValue "D29" '2.402'
$ws.Range("E29").Value = '  -0.32%  '
Set-TextValue "D30" '128.52'
$ws.Range("E30").Value = '  +2.04%  '
Set-TextValue "D31" '0.1100'
$ws.Range("E31").Value = '  +0.01%  '
Set-TextValue "D32" '1.071'
$ws.Range("E32").Value = '  -2.27%  '
Set-TextValue "D33" '3.665'
$ws.Range("E33").Value = '  -0.53%  '
Set-TextValue "D34" '5.584'
$ws.Range("E34").Value = '  -2.61%  '
Set-TextValue "D35" '0.07242'
$ws.Range("E35").Value = '  -4.84%  '
Set-TextValue "D36" '9.223'
$ws.Range("E36").Value = '  +5.26%  '
Set-TextValue "D37" '0.02344'
$ws.Range("E37").Value = '  -1.02%  '
Set-TextValue "D38" '0.2189'
$ws.Range("E38").Value = '  -1.76%  '
Set-TextValue "D39" '11.65'
$ws.Range("E39").Value = '  -3.76%  '
Set-TextValue "D40" '5.057'
$ws.Range("E40").Value = '  -3.76%  '
Set-TextValue "D41" '0.6215'
$ws.Range("E41").Value = '  -3.16%  '
Set-TextValue "D42" '1.000'
$ws.Range("E42").Value = '  -0.30%  '
Set-TextValue "D43" '1.160'
$ws.Range("E43").Value = '  -2.60%  '
Set-TextValue "D44" '13.26'
$ws.Range("E44").Value = '  -3.39%  '
Set-TextValue "D45" '0.6014'
$ws.Range("E45").Value = '  -3.62%  '
Set-TextValue "D46" '1.307'
$ws.Range("E46").Value = '  -6.71%  '
Set-TextValue "D47" '3.735'
$ws.Range("E47").Value = '  -1.66%  '
Set-TextValue "D48" '125.82'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("E49").Value = '  +1.36%  '
Set-TextValue "D50" '1.934'
$ws.Range("E50").Value = '  -3.59%  '
Set-TextValue "D51" '0.06840'
$ws.Range("E51").Value = '  -1.92%  '
